# Applies the ScheduleApproximation.xlsx edit:
#  - Inserts a new column before column D ("Working Hours") shifting the
#    "Responsible" / hours-per-person / total columns one place to the right.
#  - Inserts a new row before row 16 for a new "Testing" phase line.
#  - Updates the view selection to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at D. This shifts old columns D..K to E..L and
#    copies formatting (width, etc.) from column C, matching how Excel
#    behaves when a column is inserted via the header context menu.
$ws.Columns("D").EntireColumn.Insert()

# 2) Insert a new row at 16 for the "Testing" phase (pushes old rows
#    16-18 down to 17-19).
$ws.Rows("16").EntireRow.Insert()

$ws.Range("A16").Value = "Testing"

# Give the header cell for the new column its label.
$ws.Range("D3").Value = "Working Hours"

# 3) Restore the view state saved with the workbook.
$ws.Range("D11").Select()
